$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fix typo: "TPJ" -> "TJ" in the header row (cell I1)
$ws.Range("I1").Value = "TJ"

# Select the edited cell, matching the resulting workbook state
$ws.Range("I1").Select()
